$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IndividualBiometrics")

# Add new headers for Protein and Ontogeny columns
$ws.Range("H1").Value = "Protein"
$ws.Range("I1").Value = "Ontogeny"

# Select the newly added header cells, matching the saved selection state
$ws.Range("H1:I1").Select()

$wb.Save()
